# Insert a new data row at row 43 (pushes the old rows 43-79 down to 44-80)
# and populate it with a new "Espinaca" price record for Vega Monumental
# Concepcion, as described in the commit ("Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(43).Insert()

$ws.Range("A43").Value = 11
$ws.Range("B43").Value = "Vega Monumental Concepción"
$ws.Range("C43").Value = "Bíobío"
$ws.Range("D43").Value = 44741
$ws.Range("E43").Value = 8
$ws.Range("F43").Value = 100112012
$ws.Range("G43").Value = "Espinaca"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 100
$ws.Range("K43").Value = 8000
$ws.Range("L43").Value = 8500
$ws.Range("M43").Value = 8250
$ws.Range("N43").Value = "$/cuna 10 kilos"
$ws.Range("O43").Value = "Región Metropolitana"
$ws.Range("P43").Value = 825
$ws.Range("Q43").Value = 10
$ws.Range("R43").Value = "Hortaliza"
